# Update crypto price (D) and volume-change (E) columns per latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '64.695.34'
$c.Style = "Normal"
$ws.Range("E2").Value = '  -0.05%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.422.74'
$c.Style = "Normal"
$ws.Range("E3").Value = '  -1.05%  '

$ws.Range("E4").Value = '  -0.03%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '573.12'
$c.Style = "Normal"
$ws.Range("E5").Value = '  -0.64%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '157.12'
$c.Style = "Normal"
$ws.Range("E6").Value = '  -2.58%  '

$ws.Range("E7").Value = '  +5.93%  '

$ws.Range("E8").Value = '  -0.10%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '3.427.35'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.91%  '

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '7.16'
$c.Style = "Normal"
$ws.Range("E10").Value = '  -2.42%  '

$ws.Range("E11").Value = '  -1.66%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.443'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +0.50%  '

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '4.012.23'
$c.Style = "Normal"
$ws.Range("E13").Value = '  -1.27%  '

$ws.Range("E14").Value = '  +0.25%  '

$ws.Range("E15").Value = '  -2.87%  '

$ws.Range("E16").Value = '  -0.57%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '64.674.39'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.12%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.409.55'
$c.Style = "Normal"
$ws.Range("E18").Value = '  -1.89%  '

$ws.Range("E19").Value = '  +0.12%  '

$ws.Range("E20").Value = '  -2.05%  '

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '377.62'
$c.Style = "Normal"
$ws.Range("E21").Value = '  -2.92%  '

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '8.07'
$c.Style = "Normal"
$ws.Range("E22").Value = '  -1.25%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '0.553'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +1.15%  '

$ws.Range("E24").Value = '  +0.04%  '

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '72.45'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -0.49%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '0.0000120'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -2.87%  '

$ws.Range("E27").Value = '  +7.61%  '

$ws.Range("E28").Value = '  -1.59%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +0.02%  '

$ws.Range("E30").Value = '  +3.45%  '

$ws.Range("E31").Value = '  +0.38%  '

$ws.Range("E32").Value = '  -0.63%  '

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '23.14'
$c.Style = "Normal"
$ws.Range("E33").Value = '  -2.06%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '7.25'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +2.69%  '

$ws.Range("E35").Value = '  +7.41%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '159.99'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -1.25%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '1.91'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -0.13%  '

$ws.Range("E38").Value = '  +8.17%  '

$ws.Range("E39").Value = '  -0.74%  '

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '2.898.07'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -3.66%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '26.74'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -1.72%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '4.63'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +1.73%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '26.63'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +9.98%  '

$ws.Range("E44").Value = '  +0.23%  '

$ws.Range("E45").Value = '  -0.13%  '

$ws.Range("E46").Value = '  -0.59%  '

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '322.09'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +6.52%  '

$ws.Range("E48").Value = '  -1.42%  '

$ws.Range("E49").Value = '  +1.52%  '

$ws.Range("E50").Value = '  +1.14%  '

$ws.Range("E51").Value = '  +1.11%  '
